$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:D51 as text first so numeric-looking price strings (e.g. "594.17",
# "1.00", thousand-dot-grouped values) are stored as literal text instead of
# being auto-converted to floating point numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "70.354.78"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").Value = "3.633.24"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "594.17"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "195.54"
$ws.Range("E6").Value = "  +5.12%  "
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").Value = "3.626.16"
$ws.Range("E8").Value = "  +3.75%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("D12").Value = "58.78"
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("E14").Value = "  +4.71%  "
$ws.Range("D15").Value = "4.211.79"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("D16").Value = "19.82"
$ws.Range("E16").Value = "  +4.65%  "
$ws.Range("D17").Value = "3.632.16"
$ws.Range("E17").Value = "  +4.03%  "
$ws.Range("D18").Value = "70.360.69"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("D22").Value = "488.60"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "19.20"
$ws.Range("E23").Value = "  +15.45%  "
$ws.Range("D24").Value = "5.35"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "91.38"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("E27").Value = "  +6.80%  "
$ws.Range("D28").Value = "11.54"
$ws.Range("E28").Value = "  +4.61%  "
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  +5.24%  "
$ws.Range("D30").Value = "7.93"
$ws.Range("E30").Value = "  +10.09%  "
$ws.Range("D31").Value = "32.95"
$ws.Range("E31").Value = "  +3.97%  "
$ws.Range("E32").Value = "  +8.10%  "
$ws.Range("D33").Value = "627.07"
$ws.Range("E33").Value = "  +4.98%  "
$ws.Range("D34").Value = "12.28"
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("D35").Value = "65.97"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").Value = "40.34"
$ws.Range("E36").Value = "  +9.68%  "
$ws.Range("D37").Value = "0.414"
$ws.Range("E37").Value = "  +6.75%  "
$ws.Range("E38").Value = "  +6.99%  "
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "3.60"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "3.292.44"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").Value = "3.16"
$ws.Range("E43").Value = "  +7.93%  "
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  +11.43%  "
$ws.Range("D45").Value = "0.0453"
$ws.Range("E45").Value = "  +5.12%  "
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  +4.46%  "
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("E49").Value = "  +5.05%  "
$ws.Range("D50").Value = "3.37"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.17%  "

# Restore the plain "Normal" style on the D column so no stray number-format
# style index is left attached to the cells.
$dRange.Style = "Normal"

